$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header "Puntos de vida maximo" in column C of the header row (row 1),
# shifting the existing headers C1..AN1 one column to the right (to D1..AO1).
# Rows below the header (data row 2, formatted cell in row 7) are left untouched.
$lastCol = 40  # column AN
for ($c = $lastCol; $c -ge 3; $c--) {
    $srcCell = $ws.Cells.Item(1, $c)
    $dstCell = $ws.Cells.Item(1, $c + 1)
    $dstCell.Value = $srcCell.Value2
}

$ws.Cells.Item(1, 3).Value = "Puntos de vida maximo"

# Recompute the best-fit width for the columns whose header content changed.
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Update the selected cell shown when the workbook is reopened.
$ws.Range("F12").Select() | Out-Null
